# ACC Bid.xlsx - add "task 22" (notifikasi) and "task 23" (search) sheets
$wb = $excel.ActiveWorkbook

$task18 = $wb.Worksheets.Item(1)
$task19 = $wb.Worksheets.Item(2)

# Move the lingering selection on "task 18" from C10 to A4:B4
$task18.Range("A4:B4").Select() | Out-Null

# ---------------------------------------------------------------
# Add "task 22" sheet right after "task 19"
# ---------------------------------------------------------------
$task22 = $wb.Worksheets.Add($null, $task19)
$task22.Name = "task 22"

# Header row - reuse the bordered+bold header format from task 19 (A1:B1)
$task19.Range("A1:B1").Copy() | Out-Null
$task22.Range("A1:B1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Body rows - reuse the bordered (non-bold) body format from task 19 (A2:B2)
$task19.Range("A2:B2").Copy() | Out-Null
$task22.Range("A2:B5").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$excel.CutCopyMode = $false

$task22.Range("A2").Value = "semua notif"
$task22.Range("B2").Value = "pass"

$task22.Range("A3").Value = "hanya notif lelang"
$task22.Range("B3").Value = "pass"

$task22.Range("A4").Value = "hanya notif akun"
$task22.Range("B4").Value = "pass"

$task22.Range("A5").Value = "hanya notif aplikasi"
$task22.Range("B5").Value = "pass"

$task22.Range("A1").Value = "cek_input"
$task22.Range("B1").Value = "expected"

$task22.Columns.Item(1).ColumnWidth = 17.6
$task22.PageSetup.Orientation = 1
$task22.Range("I5").Select() | Out-Null

# ---------------------------------------------------------------
# Add "task 23" sheet right after "task 22"
# ---------------------------------------------------------------
$task23 = $wb.Worksheets.Add($null, $task22)
$task23.Name = "task 23"

$task23.Range("A1").Value = "search_text"
$task23.Range("B1").Value = "expected"
$task23.Range("A1:B1").Font.Bold = $true

$task23.Range("A2").Value = "Bagaimana cara melakukan Topup?"
$task23.Range("B2").Value = "pass"

$task23.Range("B3").Value = "fail"

$task23.Range("A4").Value = "xx92923"
$task23.Range("B4").Value = "fail"

$task23.Columns.Item(1).ColumnWidth = 31.6
$task23.Range("A1:B4").Select() | Out-Null

# Make "task 23" the active sheet/tab (last-saved view state)
$task23.Activate()
